{"js": "// 1) Remove the stray \"_GoBack\" bookmark that currently sits at the end of\n//    the first paragraph (after \"SETI Institute\").\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Remove the justified alignment from the abstract paragraph (2nd\n//    paragraph) -- it goes back to the (unset/left) default alignment.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst abstractParagraph = paragraphs.items[1];\nabstractParagraph.alignment = Word.Alignment.left;\nawait context.sync();\n\n// 3) Re-insert the \"_GoBack\" bookmark at its new location: right after the\n//    \"g\" in \"preserving\" (i.e. spanning just that \"g\"), splitting the run\n//    that used to read \"...while preserving astrophysical...\".\nconst searchResults = context.document.body.search(\"while preservin\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nconst hit = searchResults.items[0];\nconst afterHit = hit.getRange(\"After\");\nconst gRange = afterHit.getRange(\"Start\").expandTo(afterHit.getRange(\"Start\"));\n// Select exactly the \"g\" character right after the found text, then bookmark it.\nconst gOnly = afterHit.insertText(\"g\", \"Before\");\ngOnly.insertBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# 1) Remove the stray \"_GoBack\" bookmark that currently sits at the end of\n#    the first paragraph (after \"SETI Institute\").\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n# 2) Remove the justified alignment from the abstract paragraph (2nd\n#    paragraph) -- it goes back to left/default alignment.\n$d.Paragraphs(2).Alignment = 0  # wdAlignParagraphLeft\n\n# 3) Re-insert the \"_GoBack\" bookmark at its new location: right after the\n#    \"g\" in \"preserving\" (i.e. spanning just that \"g\"), inside the run that\n#    used to read \"...while preserving astrophysical...\".\n$rng = $d.Content\n$find = $rng.Find\n$find.ClearFormatting()\n$find.Text = \"while preservin\"\n$find.Forward = $true\n$find.Wrap = 1  # wdFindContinue\n$found = $find.Execute()\nif ($found) {\n    $markStart = $rng.End\n    $gRange = $d.Range($markStart, $markStart + 1)\n    $d.Bookmarks.Add(\"_GoBack\", $gRange)\n}\n"}
